$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version, Date, Contact values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Sheet "Include from EntityNameUse": insert a new "ASGN" concept row after "SNDX" (row 10) ---
$incl = $wb.Worksheets.Item("Include from EntityNameUse")
$incl.Rows.Item(11).Insert()

# Copy formatting from the row above (SNDX row) so the new row matches the existing style
$incl.Range("A10:B10").Copy()
$incl.Range("A11:B11").PasteSpecial(-4122)

$incl.Range("A11").Value = "ASGN"

Write-Output "done"
